# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-column suffixes to "_FV2210" / "_FV2304"
# - Turn the header range A1:U75 into an Excel Table ("Table1") with AutoFilter
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row -------------------------------------------------
# Columns A:J were "<Label>_old"  -> "<Label>_FV2210"
# Column  K stays "diff"
# Columns L:U were "<Label>_new"  -> "<Label>_FV2304"

$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $label = $cell.Value()
    $cell.Value = ($label -replace "_old$", "_FV2210")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $label = $cell.Value()
    $cell.Value = ($label -replace "_new$", "_FV2304")
}

# --- 2. Turn A1:U75 into a real Table (ListObject) named "Table1" ------------

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U75"), $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (row 1) -----------------------------------------

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
